$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new "season record" columns, matching the
# style used by the existing header row (same bold/border/center style
# as the rest of row 1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill the season record for every player row (2-40) with the team's
# Wins / Losses / Ties totals.
$lastRow = 40
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 86
    $ws.Cells.Item($r, 31).Value = 76
    $ws.Cells.Item($r, 32).Value = 0
}
